# "First Final Draft of chapters 1 and 2"
#
# Slide 2 ("What are extras?"): the "extras?" word moves out of the Title
# placeholder (where it sat on its own line) and becomes the caption of the
# "Oval 4" bubble shape, which is repositioned to sit beside/over the title.
# The title shape itself gains an explicit position/size, and shape z-order
# changes so Title, then Oval 4, then Oval 5 are the first three shapes.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

$title = $s.Shapes.Item("Title 1")
$oval4 = $s.Shapes.Item("Oval 4")

# --- Title 1: give it an explicit position/size, and drop the manual line
#     break + "extras?" run so it just reads "What are" ---
$title.Left = 36
$title.Top = 18
$title.Width = 504
$title.Height = 90

$trTitle = $title.TextFrame.TextRange
$fullTitle = $trTitle.Characters(1, $trTitle.Length)
$fullTitle.Text = "What are"

# --- Oval 4: move it next to the title, and give it the "extras?" caption ---
$oval4.Left = 378
$oval4.Top = 30

$trOval4 = $oval4.TextFrame.TextRange
$trOval4.Text = "extras?"
$trOval4.Font.Size = 44

# --- Re-order shapes: Title, Oval 4, Oval 5, Content Placeholder, Picture ---
$title.ZOrder(1)   # msoSendToBack -> Title becomes shape 1
$oval4.ZOrder(3)   # msoSendBackward (twice) -> Oval 4 becomes shape 2,
$oval4.ZOrder(3)   #   right after the Title and ahead of Oval 5
